$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27:89 down to 28:90
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record
$ws.Cells.Item(27, 1).Value = 11
$ws.Cells.Item(27, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(27, 3).Value = "Bíobío"
$ws.Cells.Item(27, 4).Value = "2022-05-27"
$ws.Cells.Item(27, 5).Value = 8
$ws.Cells.Item(27, 6).Value = 100112001
$ws.Cells.Item(27, 7).Value = "Berenjena"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 50
$ws.Cells.Item(27, 11).Value = 5000
$ws.Cells.Item(27, 12).Value = 5500
$ws.Cells.Item(27, 13).Value = 5200
$ws.Cells.Item(27, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 87
$ws.Cells.Item(27, 17).Value = 60
$ws.Cells.Item(27, 18).Value = "Hortaliza"
